# Applies a row-permutation of columns A,B,D,E,F,G,H,Q,R across rows 127-161
# (row 158 is left untouched). For each target row, the new values are taken
# from the ORIGINAL (pre-edit) values of the row given by $mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row's ORIGINAL values get copied into target row)
$mapping = @{
    127 = 142
    128 = 140
    129 = 152
    130 = 151
    131 = 154
    132 = 143
    133 = 127
    134 = 153
    135 = 155
    136 = 156
    137 = 129
    138 = 146
    139 = 160
    140 = 135
    141 = 136
    142 = 137
    143 = 157
    144 = 138
    145 = 150
    146 = 149
    147 = 141
    148 = 132
    149 = 161
    150 = 139
    151 = 134
    152 = 128
    153 = 144
    154 = 131
    155 = 145
    156 = 159
    157 = 130
    159 = 147
    160 = 148
    161 = 133
}

$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)  # A, B, D, E, F, G, H, Q, R

# Snapshot the original values for every row involved (both as source and target)
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Now write the permuted values into every target row, using ONLY the snapshot
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcVals[$c]
    }
}
